$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 11.8
$ws.Range("O2").Value = 1.31
$ws.Range("P2").Value = 3.47
$ws.Range("U2").Value = 1.75
$ws.Range("V2").Value = 1.91
